$wb = $excel.ActiveWorkbook

# --- Sheet "L_map" : update a couple of latex-name labels (V_R -> V_LPF family) ---
$wsLmap = $wb.Worksheets.Item("L_map")
[void]$wsLmap.Activate()
$wsLmap.Range("E5").Value = "V_LPF"
$wsLmap.Range("F5").Value = "V_{LPF}"
$wsLmap.Range("E15").Value = "V_LPF"

# --- Sheet "var" : bump the w_c value, leave selection on B8 ---
$wsVar = $wb.Worksheets.Item("var")
[void]$wsVar.Activate()
$wsVar.Range("B7").Value = 2800
[void]$wsVar.Range("B8").Select()

# --- Sheet "SS0" : move the selection only ---
$wsSS0 = $wb.Worksheets.Item("SS0")
[void]$wsSS0.Activate()
[void]$wsSS0.Range("C5").Select()

# --- Sheet "SS1" : no content/selection change in the target state ---

# --- Sheet "SS2" : no content/selection change in the target state ---

# --- Finish back on "L_map", matching the final active sheet/selection ---
[void]$wsLmap.Activate()
[void]$wsLmap.Range("F4").Select()
